$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1): new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, thin borders, center/top alignment)
# from the existing H1 header cell onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data cells (row 2): plain numeric values, no explicit style (matches C2..H2)
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
